$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 78
$ws.Range("I2").Value = 177
$ws.Range("J2").Value = 755
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 207
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = 132
$ws.Range("O2").Value = 2
$ws.Range("R2").Value = 11
$ws.Range("S2").Value = 62
$ws.Range("T2").Value = 118
$ws.Range("V2").Value = 1094
$ws.Range("X2").Value = 1098
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 10
$ws.Range("AA2").Value = 11
